# Weekly data refresh: insert the newest week's row at the top of the
# data block (row 79), pushing all existing rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 79 (existing rows 79..113
# shift down to 80..114; formats/styles are inherited from the row below).
$ws.Rows.Item(79).Insert()

# Populate the new row 79 with this week's record.
$ws.Range("A79").Value = 1
$ws.Range("B79").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C79").Value = "Arica y Parinacota"
$ws.Range("D79").Value = 44680
$ws.Range("E79").Value = 15
$ws.Range("F79").Value = "Fruta"
$ws.Range("G79").Value = 100108
$ws.Range("H79").Value = "Tropicales y subtropicales"
$ws.Range("I79").Value = 100108003
$ws.Range("J79").Value = "Maracuyá"
$ws.Range("K79").Value = "Sin especificar"
$ws.Range("L79").Value = "Primera"
$ws.Range("M79").Value = 120
$ws.Range("N79").Value = 21000
$ws.Range("O79").Value = 22000
$ws.Range("P79").Value = 21500
$ws.Range("Q79").Value = "$/caja 20 kilos"
$ws.Range("R79").Value = "Región de Arica y Parinacota"
$ws.Range("S79").Value = 1075
$ws.Range("T79").Value = 20
